# Slide 7 (sldId 647) - "Rectangle 3" text box (shape 4):
# Split the single-run "#define BASE_FARE 3.20" and
# "#define SQUARE(x) x*x" lines into three colour-coded runs each,
# matching the "#define NAME value" syntax-highlighting style already
# used elsewhere on the slide (green keyword / blue identifier / red value).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$tr = $s.Shapes.Item(4).TextFrame.TextRange

# Paragraph 4: "#define BASE_FARE 3.20"
$para4 = $tr.Paragraphs(4, 1)
$para4.Characters(1, 8).Font.Color.RGB = 26112      # "#define "  -> 006600 (green)
$para4.Characters(9, 10).Font.Color.RGB = 16711680  # "BASE_FARE " -> 0000FF (blue)
$para4.Characters(19, 4).Font.Color.RGB = 192       # "3.20"      -> C00000 (red, unchanged)

# Paragraph 5: "#define SQUARE(x) x*x"
$para5 = $tr.Paragraphs(5, 1)
$para5.Characters(1, 8).Font.Color.RGB = 26112      # "#define "    -> 006600 (green)
$para5.Characters(9, 10).Font.Color.RGB = 16711680  # "SQUARE(x) " -> 0000FF (blue)
$para5.Characters(19, 3).Font.Color.RGB = 192       # "x*x"        -> C00000 (red, unchanged)
